$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# sheet1 changes
$ws1.Range('A2').Value = 'DTaP/'
$ws1.Range('A3').Value = 'DTaP/'
$ws1.Range('A4').Value = 'DTaP/'
$ws1.Range('A5').Value = 'DTaP/'
$ws1.Range('A8').Value = 'DTaP-Hib '
$ws1.Range('A11').Value = 'Hepatitis B-Hib'
$ws1.Range('A15').Value = 'Hepatitis A-Hepatitis B 18 only'
$ws1.Range('A16').Value = 'Hepatitis A-Hepatitis B 18 only'
$ws1.Range('A17').Value = 'Hepatitis B Pediatric/Adolescent'
$ws1.Range('A18').Value = 'Hepatitis B Pediatric/Adolescent'
$ws1.Range('A19').Value = 'Hepatitis B Pediatric/Adolescent'
$ws1.Range('A25').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws1.Range('A26').Value = 'Meningococcal Conjugate (Groups A, C, Y and W-135) '
$ws1.Range('A27').Value = 'Measles, Mumps and Rubella (MMR)/'
$ws1.Range('A28').Value = 'Pneumococcal 7-valent (Pediatric)'
$ws1.Range('A30').Value = 'Tetanus  Diphtheria Toxoids'
$ws1.Range('D30').Value = '10 pack - 1 dose syringes No Needle '
$ws1.Range('A31').Value = 'Tetanus  Diphtheria Toxoids'
$ws1.Range('D31').Value = '10 pack - 1 dose vials '
$ws1.Range('A32').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/'
$ws1.Range('D32').Value = '10 pack - 1 dose vials '
$ws1.Range('A33').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/'
$ws1.Range('D33').Value = '5 pack - 1 dose TL syringes, No Needle '
$ws1.Range('A34').Value = 'Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/'

# sheet2 changes
$ws2.Range('A6').Value = 'Hepatitis A-Hepatitis B Adult'
$ws2.Range('A7').Value = 'Hepatitis A-Hepatitis B Adult'
$ws2.Range('A13').Value = 'Tetanus  Diphtheria Toxoids'

# sheet3 changes
$ws3.Range('A2').Value = 'Influenza   (Age 6 months and older)'
$ws3.Range('A3').Value = 'Influenza  (Age 6-35 months)'
$ws3.Range('A4').Value = 'Influenza  (Age 36 months and older)'
$ws3.Range('A5').Value = 'Influenza  (Age 36 months and older)'
$ws3.Range('A6').Value = 'Influenza   (Age 4 years and older)'
$ws3.Range('A7').Value = 'Influenza  (Age 18 years and older)'
$ws3.Range('A8').Value = 'Influenza  Live, Intranasal (Age 5-49 years)'
